$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12490.173
$ws.Range("J17").Value = 12642.743
$ws.Range("L17").Value = 37928.229
$ws.Range("N17").Value = -38264.229
$ws.Range("H64").Value = 3630.7307
$ws.Range("I64").Value = 3535.96
$ws.Range("K64").Value = 3535.96
$ws.Range("M64").Value = -3287.96
$ws.Range("H67").Value = 3630.7307
$ws.Range("I67").Value = 3535.96
$ws.Range("K67").Value = 3535.96
$ws.Range("M67").Value = -2677.96
$ws.Range("H76").Value = 3991.5833
$ws.Range("I76").Value = 3440.6
$ws.Range("J76").Value = 6746.5
$ws.Range("K76").Value = 3440.6
$ws.Range("L76").Value = 6746.5
$ws.Range("M76").Value = -3125.6
$ws.Range("N76").Value = -7376.5
$ws.Range("H79").Value = 3991.5833
$ws.Range("I79").Value = 3440.6
$ws.Range("J79").Value = 6746.5
$ws.Range("K79").Value = 3440.6
$ws.Range("L79").Value = 6746.5
$ws.Range("M79").Value = -2348.6
$ws.Range("N79").Value = -8930.5
$ws.Range("H86").Value = 2689.7
$ws.Range("I86").Value = 2685.2856
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 2685.2856
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -1562.2856
$ws.Range("N86").Value = -4946
$ws.Range("H89").Value = 2689.7
$ws.Range("I89").Value = 2685.2856
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 13426.428
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -7810.428
$ws.Range("N89").Value = -24732
$ws.Range("H92").Value = 1100.25
$ws.Range("I92").Value = 1114.5714
$ws.Range("K92").Value = 1114.5714
$ws.Range("M92").Value = 133.4286
$ws.Range("H98").Value = 510.6
$ws.Range("I98").Value = 511.1579
$ws.Range("K98").Value = 511.1579
$ws.Range("M98").Value = 986.8421000000001
$ws.Range("H106").Value = 16538.715
$ws.Range("I106").Value = 1925.6666
$ws.Range("J106").Value = 27498.5
$ws.Range("K106").Value = 1925.6666
$ws.Range("L106").Value = 27498.5
$ws.Range("M106").Value = -1294.6666
$ws.Range("N106").Value = -28760.5
$ws.Range("H112").Value = 2321.6667
$ws.Range("I112").Value = 1408.3334
$ws.Range("J112").Value = 2452.1428
$ws.Range("K112").Value = 4225.0002
$ws.Range("L112").Value = 7356.428400000001
$ws.Range("M112").Value = -3117.0002
$ws.Range("N112").Value = -9572.428400000001
$ws.Range("H116").Value = 5375
$ws.Range("I116").Value = 7500
$ws.Range("J116").Value = 4950
$ws.Range("K116").Value = 7500
$ws.Range("L116").Value = 4950
$ws.Range("M116").Value = -4058
$ws.Range("N116").Value = -11834
$ws.Range("H122").Value = 510.6
$ws.Range("I122").Value = 511.1579
$ws.Range("K122").Value = 1533.4737
$ws.Range("M122").Value = 916.5263
$ws.Range("H127").Value = 102388.78
$ws.Range("I127").Value = 102388.78
$ws.Range("K127").Value = 307166.34
$ws.Range("M127").Value = -302206.34
$ws.Range("H132").Value = 3347.75
$ws.Range("I132").Value = 1298.2142
$ws.Range("J132").Value = 8130
$ws.Range("K132").Value = 3894.6426
$ws.Range("L132").Value = 24390
$ws.Range("M132").Value = -1364.6426
$ws.Range("N132").Value = -29450
$ws.Range("H137").Value = 3172.6123
$ws.Range("I137").Value = 1152
$ws.Range("K137").Value = 3456
$ws.Range("M137").Value = -906
$ws.Range("H138").Value = 2005.3103
$ws.Range("I138").Value = 1310.8667
$ws.Range("K138").Value = 3932.6001
$ws.Range("M138").Value = 1207.3999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35983.67
$ws.Range("I32").Value = 20444.703
$ws.Range("K32").Value = 20444.703
$ws.Range("M32").Value = -20157.703
$ws.Range("H74").Value = 1727.5897
$ws.Range("I74").Value = 812.73334
$ws.Range("K74").Value = 812.73334
$ws.Range("M74").Value = 61.26666
$ws.Range("H77").Value = 1727.5897
$ws.Range("I77").Value = 812.73334
$ws.Range("K77").Value = 4063.6667
$ws.Range("M77").Value = 304.3333000000002
$ws.Range("H122").Value = 3586.9092
$ws.Range("I122").Value = 2307.125
$ws.Range("J122").Value = 6999.6665
$ws.Range("K122").Value = 6921.375
$ws.Range("L122").Value = 20998.9995
$ws.Range("M122").Value = -4471.375
$ws.Range("N122").Value = -25898.9995
$ws.Range("H132").Value = 13660.412
$ws.Range("I132").Value = 16632.77
$ws.Range("J132").Value = 4000.25
$ws.Range("K132").Value = 49898.31
$ws.Range("L132").Value = 12000.75
$ws.Range("M132").Value = -47368.31
$ws.Range("N132").Value = -17060.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10958.6875
$ws.Range("I20").Value = 7638.25
$ws.Range("J20").Value = 14279.125
$ws.Range("K20").Value = 7638.25
$ws.Range("L20").Value = 14279.125
$ws.Range("M20").Value = -7391.25
$ws.Range("N20").Value = -14773.125
$ws.Range("H35").Value = 27142.857
$ws.Range("J35").Value = 27142.857
$ws.Range("L35").Value = 27142.857
$ws.Range("N35").Value = -27762.857
$ws.Range("H86").Value = 3945.5715
$ws.Range("I86").Value = 3963.4
$ws.Range("J86").Value = 3901
$ws.Range("K86").Value = 3963.4
$ws.Range("L86").Value = 3901
$ws.Range("M86").Value = -2840.4
$ws.Range("N86").Value = -6147
$ws.Range("H89").Value = 3945.5715
$ws.Range("I89").Value = 3963.4
$ws.Range("J89").Value = 3901
$ws.Range("K89").Value = 19817
$ws.Range("L89").Value = 19505
$ws.Range("M89").Value = -14201
$ws.Range("N89").Value = -30737
$ws.Range("H94").Value = 657.3077
$ws.Range("I94").Value = 635.4211
$ws.Range("J94").Value = 716.7143
$ws.Range("K94").Value = 635.4211
$ws.Range("L94").Value = 716.7143
$ws.Range("M94").Value = -184.4211
$ws.Range("N94").Value = -1618.7143
$ws.Range("H107").Value = 20205.537
$ws.Range("I107").Value = 33121.72
$ws.Range("J107").Value = 1418.3636
$ws.Range("K107").Value = 33121.72
$ws.Range("L107").Value = 1418.3636
$ws.Range("M107").Value = -31201.72
$ws.Range("N107").Value = -5258.3636
$ws.Range("H134").Value = 1402.625
$ws.Range("I134").Value = 1388.762
$ws.Range("K134").Value = 4166.286
$ws.Range("M134").Value = -1631.286

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2878.5
$ws.Range("I16").Value = 3098.25
$ws.Range("K16").Value = 3098.25
$ws.Range("M16").Value = -2811.25
$ws.Range("H31").Value = 1418.7241
$ws.Range("I31").Value = 1270.5
$ws.Range("J31").Value = 2130.2
$ws.Range("K31").Value = 1270.5
$ws.Range("L31").Value = 2130.2
$ws.Range("M31").Value = -975.5
$ws.Range("N31").Value = -2720.2
$ws.Range("H34").Value = 1418.7241
$ws.Range("I34").Value = 1270.5
$ws.Range("J34").Value = 2130.2
$ws.Range("K34").Value = 1270.5
$ws.Range("L34").Value = 2130.2
$ws.Range("M34").Value = -1068.5
$ws.Range("N34").Value = -2534.2
$ws.Range("H58").Value = 1141.4166
$ws.Range("I58").Value = 898.5
$ws.Range("J58").Value = 1627.25
$ws.Range("K58").Value = 898.5
$ws.Range("L58").Value = 1627.25
$ws.Range("M58").Value = -695.5
$ws.Range("N58").Value = -2033.25
$ws.Range("H62").Value = 3788.111
$ws.Range("I62").Value = 3813.4285
$ws.Range("K62").Value = 3813.4285
$ws.Range("M62").Value = -3189.4285
$ws.Range("H65").Value = 3788.111
$ws.Range("I65").Value = 3813.4285
$ws.Range("K65").Value = 19067.1425
$ws.Range("M65").Value = -15947.1425
$ws.Range("H86").Value = 6707.6924
$ws.Range("I86").Value = 5633.4443
$ws.Range("J86").Value = 9124.75
$ws.Range("K86").Value = 5633.4443
$ws.Range("L86").Value = 9124.75
$ws.Range("M86").Value = -4510.4443
$ws.Range("N86").Value = -11370.75
$ws.Range("H89").Value = 6707.6924
$ws.Range("I89").Value = 5633.4443
$ws.Range("J89").Value = 9124.75
$ws.Range("K89").Value = 28167.2215
$ws.Range("L89").Value = 45623.75
$ws.Range("M89").Value = -22551.2215
$ws.Range("N89").Value = -56855.75
$ws.Range("H99").Value = 2119.0952
$ws.Range("I99").Value = 1590.4
$ws.Range("J99").Value = 2599.7273
$ws.Range("K99").Value = 1590.4
$ws.Range("L99").Value = 2599.7273
$ws.Range("M99").Value = -92.40000000000009
$ws.Range("N99").Value = -5595.7273
$ws.Range("H105").Value = 1803.8
$ws.Range("I105").Value = 1650
$ws.Range("K105").Value = 1650
$ws.Range("M105").Value = 97
$ws.Range("H113").Value = 2878.5
$ws.Range("I113").Value = 3098.25
$ws.Range("K113").Value = 3098.25
$ws.Range("M113").Value = -928.25
$ws.Range("H126").Value = 2119.0952
$ws.Range("I126").Value = 1590.4
$ws.Range("J126").Value = 2599.7273
$ws.Range("K126").Value = 4771.200000000001
$ws.Range("L126").Value = 7799.1819
$ws.Range("M126").Value = -2301.200000000001
$ws.Range("N126").Value = -12739.1819
$ws.Range("H132").Value = 3601.111
$ws.Range("I132").Value = 3675.8333
$ws.Range("K132").Value = 11027.4999
$ws.Range("M132").Value = -8497.499899999999
$ws.Range("H134").Value = 2563.9167
$ws.Range("I134").Value = 2342.4546
$ws.Range("K134").Value = 7027.3638
$ws.Range("M134").Value = -4492.3638
$ws.Range("H136").Value = 1141.4166
$ws.Range("I136").Value = 898.5
$ws.Range("J136").Value = 1627.25
$ws.Range("K136").Value = 2695.5
$ws.Range("L136").Value = 4881.75
$ws.Range("M136").Value = -145.5
$ws.Range("N136").Value = -9981.75
$ws.Range("H141").Value = 107553.2
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 9924.6
$ws.Range("I25").Value = 312
$ws.Range("J25").Value = 16333
$ws.Range("K25").Value = 936
$ws.Range("L25").Value = 48999
$ws.Range("M25").Value = -767
$ws.Range("N25").Value = -49337
$ws.Range("H30").Value = 9924.6
$ws.Range("I30").Value = 312
$ws.Range("J30").Value = 16333
$ws.Range("K30").Value = 936
$ws.Range("L30").Value = 48999
$ws.Range("M30").Value = -834
$ws.Range("N30").Value = -49203
$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 4000
$ws.Range("K63").Value = 12000
$ws.Range("M63").Value = -11251
$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 4000
$ws.Range("K66").Value = 36000
$ws.Range("M66").Value = -32256
$ws.Range("H75").Value = 4539.4
$ws.Range("I75").Value = 3301
$ws.Range("J75").Value = 4849
$ws.Range("K75").Value = 9903
$ws.Range("L75").Value = 14547
$ws.Range("M75").Value = -8905
$ws.Range("N75").Value = -16543
$ws.Range("H78").Value = 4539.4
$ws.Range("I78").Value = 3301
$ws.Range("J78").Value = 4849
$ws.Range("K78").Value = 29709
$ws.Range("L78").Value = 43641
$ws.Range("M78").Value = -24717
$ws.Range("N78").Value = -53625
$ws.Range("H92").Value = 2896.7144
$ws.Range("J92").Value = 2894.8
$ws.Range("L92").Value = 8684.400000000001
$ws.Range("N92").Value = -11180.4
$ws.Range("H129").Value = 2644.6316
$ws.Range("I129").Value = 572.1111
$ws.Range("J129").Value = 4509.9
$ws.Range("K129").Value = 1716.3333
$ws.Range("L129").Value = 13529.7
$ws.Range("M129").Value = 3283.6667
$ws.Range("N129").Value = -23529.7
$ws.Range("H131").Value = 2241.6667
$ws.Range("J131").Value = 2450
$ws.Range("L131").Value = 7350
$ws.Range("N131").Value = -17430
$ws.Range("H132").Value = 1668.762
$ws.Range("I132").Value = 759
$ws.Range("J132").Value = 2351.0833
$ws.Range("K132").Value = 6831
$ws.Range("L132").Value = 21159.7497
$ws.Range("M132").Value = -4301
$ws.Range("N132").Value = -26219.7497
$ws.Range("H137").Value = 3678.9
$ws.Range("I137").Value = 1578.2858
$ws.Range("J137").Value = 4810
$ws.Range("K137").Value = 4734.857400000001
$ws.Range("L137").Value = 14430
$ws.Range("M137").Value = 365.1425999999992
$ws.Range("N137").Value = -24630
$ws.Range("H138").Value = 2332.75
$ws.Range("I138").Value = 1325.8
$ws.Range("K138").Value = 3977.4
$ws.Range("M138").Value = 1162.6
$ws.Range("H139").Value = 2198.2
$ws.Range("I139").Value = 1997.9231
$ws.Range("K139").Value = 5993.7693
$ws.Range("M139").Value = -853.7692999999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1658.3334
$ws.Range("I97").Value = 985.8333
$ws.Range("J97").Value = 3003.3333
$ws.Range("K97").Value = 985.8333
$ws.Range("L97").Value = 3003.3333
$ws.Range("M97").Value = -489.8333
$ws.Range("N97").Value = -3995.3333
$ws.Range("H107").Value = 20835064
$ws.Range("I107").Value = 529.3
$ws.Range("J107").Value = 35716876
$ws.Range("K107").Value = 529.3
$ws.Range("L107").Value = 35716876
$ws.Range("M107").Value = 1390.7
$ws.Range("N107").Value = -35720716
$ws.Range("H113").Value = 4305.5
$ws.Range("I113").Value = 3474.25
$ws.Range("K113").Value = 3474.25
$ws.Range("M113").Value = -1304.25
$ws.Range("H122").Value = 194154.84
$ws.Range("I122").Value = 239705.28
$ws.Range("J122").Value = 2843
$ws.Range("K122").Value = 719115.84
$ws.Range("L122").Value = 8529
$ws.Range("M122").Value = -716665.84
$ws.Range("N122").Value = -13429
$ws.Range("H132").Value = 2766.4
$ws.Range("I132").Value = 2558.4443
$ws.Range("K132").Value = 7675.3329
$ws.Range("M132").Value = -5145.3329

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5011.3687
$ws.Range("I7").Value = 5587.5454
$ws.Range("J7").Value = 4219.125
$ws.Range("K7").Value = 5587.5454
$ws.Range("L7").Value = 4219.125
$ws.Range("M7").Value = -5475.5454
$ws.Range("N7").Value = -4443.125
$ws.Range("H40").Value = 3769.2942
$ws.Range("I40").Value = 2606.5
$ws.Range("J40").Value = 6560
$ws.Range("K40").Value = 2606.5
$ws.Range("L40").Value = 6560
$ws.Range("M40").Value = -2470.5
$ws.Range("N40").Value = -6832
$ws.Range("H46").Value = 22637.85
$ws.Range("I46").Value = 31768.785
$ws.Range("J46").Value = 1332.3334
$ws.Range("K46").Value = 31768.785
$ws.Range("L46").Value = 1332.3334
$ws.Range("M46").Value = -31580.785
$ws.Range("N46").Value = -1708.3334
$ws.Range("H55").Value = 857.5
$ws.Range("I55").Value = 636.6667
$ws.Range("J55").Value = 1078.3334
$ws.Range("K55").Value = 636.6667
$ws.Range("L55").Value = 1078.3334
$ws.Range("M55").Value = -463.6667
$ws.Range("N55").Value = -1424.3334
$ws.Range("H68").Value = 2111.077
$ws.Range("I68").Value = 1817
$ws.Range("K68").Value = 1817
$ws.Range("M68").Value = -1068
$ws.Range("H71").Value = 2111.077
$ws.Range("I71").Value = 1817
$ws.Range("K71").Value = 9085
$ws.Range("M71").Value = -5341
$ws.Range("H93").Value = 21684.412
$ws.Range("I93").Value = 2221.2856
$ws.Range("J93").Value = 112512.336
$ws.Range("K93").Value = 2221.2856
$ws.Range("L93").Value = 112512.336
$ws.Range("M93").Value = -973.2856000000002
$ws.Range("N93").Value = -115008.336
$ws.Range("H122").Value = 3081
$ws.Range("I122").Value = 2628.7
$ws.Range("J122").Value = 3834.8333
$ws.Range("K122").Value = 7886.099999999999
$ws.Range("L122").Value = 11504.4999
$ws.Range("M122").Value = -5436.099999999999
$ws.Range("N122").Value = -16404.4999
$ws.Range("H126").Value = 5011.3687
$ws.Range("I126").Value = 5587.5454
$ws.Range("J126").Value = 4219.125
$ws.Range("K126").Value = 16762.6362
$ws.Range("L126").Value = 12657.375
$ws.Range("M126").Value = -14292.6362
$ws.Range("N126").Value = -17597.375
$ws.Range("H136").Value = 2419.75
$ws.Range("I136").Value = 1870.5333
$ws.Range("J136").Value = 5165.8335
$ws.Range("K136").Value = 5611.5999
$ws.Range("L136").Value = 15497.5005
$ws.Range("M136").Value = -3061.5999
$ws.Range("N136").Value = -20597.5005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3346.0557
$ws.Range("I81").Value = 3748.6
$ws.Range("J81").Value = 1333.3334
$ws.Range("K81").Value = 7497.2
$ws.Range("L81").Value = 2666.6668
$ws.Range("M81").Value = -6436.2
$ws.Range("N81").Value = -4788.6668
$ws.Range("H84").Value = 3346.0557
$ws.Range("I84").Value = 3748.6
$ws.Range("J84").Value = 1333.3334
$ws.Range("K84").Value = 37486
$ws.Range("L84").Value = 13333.334
$ws.Range("M84").Value = -32182
$ws.Range("N84").Value = -23941.334
$ws.Range("H92").Value = 33439.8
$ws.Range("J92").Value = 33439.8
$ws.Range("L92").Value = 33439.8
$ws.Range("N92").Value = -38431.8
$ws.Range("H96").Value = 3703.6155
$ws.Range("I96").Value = 3679.0833
$ws.Range("J96").Value = 3998
$ws.Range("K96").Value = 3679.0833
$ws.Range("L96").Value = 3998
$ws.Range("M96").Value = -2306.0833
$ws.Range("N96").Value = -6744
$ws.Range("H107").Value = 31250802
$ws.Range("I107").Value = 820.2857
$ws.Range("J107").Value = 55556344
$ws.Range("K107").Value = 2460.8571
$ws.Range("L107").Value = 166669032
$ws.Range("M107").Value = -540.8571000000002
$ws.Range("N107").Value = -166672872
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 2455.1667
$ws.Range("I122").Value = 2033.4286
$ws.Range("J122").Value = 3931.25
$ws.Range("K122").Value = 6100.2858
$ws.Range("L122").Value = 11793.75
$ws.Range("M122").Value = -3650.2858
$ws.Range("N122").Value = -16693.75
$ws.Range("H132").Value = 17367.71
$ws.Range("I132").Value = 21751.834
$ws.Range("K132").Value = 65255.50199999999
$ws.Range("M132").Value = -62725.50199999999
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws.Range("H136").Value = 2412.1
$ws.Range("I136").Value = 2559.3928
$ws.Range("J136").Value = 350
$ws.Range("K136").Value = 7678.178400000001
$ws.Range("L136").Value = 1050
$ws.Range("M136").Value = -5128.178400000001
$ws.Range("N136").Value = -6150
